# Auto-committed on 2022/03/14 (週一)
# Fills in the review-tracking sheet for rows 75-83 ("User審查意見彙整"):
#   F = SKL's reply/explanation text for each comment
#   G = expected-completion date (預計完成日)
#   H = reply date (意見回覆日)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User審查意見彙整")

# --- F column: reply / explanation text ---------------------------------
$ws.Range("F75").Value = "已修改"
$ws.Range("F76").Value = "增加查詢按鈕供查看資料"
$ws.Range("F77").Value = "已修改"
$ws.Range("F78").Value = "UAT環境，可查詢入帳日期=1101229"
$ws.Range("F79").Value = "僅供依此條件查看延遲交易確認資料，若需移除條件也可以"
$ws.Range("F80").Value = "原為因實際還款日有值(自動帶入)則金額必須輸入，現實際還款日不自動帶入"
$ws.Range("F81").Value = "原為預設日期，可不輸入，已經預設部分移除"
$ws.Range("F82").Value = "已修改"
$ws.Range("F83").Value = "已修改"

# --- G column: 預計完成日 (only row 75 was still blank) -------------------
$ws.Range("G75").Value = "2022/3/16"

# --- H column: 意見回覆日 for the whole block -----------------------------
$ws.Range("H75").Value = "2022/3/14"
$ws.Range("H76").Value = "2022/3/14"
$ws.Range("H77").Value = "2022/3/14"
$ws.Range("H78").Value = "2022/3/14"
$ws.Range("H79").Value = "2022/3/14"
$ws.Range("H80").Value = "2022/3/14"
$ws.Range("H81").Value = "2022/3/14"
$ws.Range("H82").Value = "2022/3/14"
$ws.Range("H83").Value = "2022/3/14"

# --- Update the on-screen view to where the author was working ----------
$ws.Activate()
$ws.Range("H78").Select()
